$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --- Cells changing from numeric to the shared "0" text marker ---
# (copies format+value from a donor cell that already holds text "0" with the correct style)
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("C15").Copy($ws.Range("C18"))
$ws.Range("C15").Copy($ws.Range("C27"))

# --- Cell changing from numeric to the shared "***.*" text marker ---
$ws.Range("E14").Copy($ws.Range("E15"))

# --- Cell reverting from the "0" text marker back to a plain number ---
$ws.Range("D23").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 3

# --- Plain numeric value updates ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 22.222222222222
$ws.Range("I16").Value = 134
$ws.Range("J16").Value = 164
$ws.Range("K16").Value = -18.292682926829
$ws.Range("L16").Value = 47.252747252747
$ws.Range("M16").Value = -8.843537414965
$ws.Range("N16").Value = -75.092936802974
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -6.666666666666
$ws.Range("I17").Value = 180
$ws.Range("J17").Value = 189
$ws.Range("K17").Value = -4.761904761904
$ws.Range("L17").Value = 29.496402877697
$ws.Range("M17").Value = 69.811320754717
$ws.Range("N17").Value = -22.746781115879
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 117
$ws.Range("K18").Value = 6.837606837606
$ws.Range("L18").Value = 56.25
$ws.Range("M18").Value = -13.793103448275
$ws.Range("N18").Value = -87.474949899799
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 77.777777777777
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = 5.357142857142
$ws.Range("I19").Value = 589
$ws.Range("J19").Value = 607
$ws.Range("K19").Value = -2.965403624382
$ws.Range("L19").Value = 45.073891625615
$ws.Range("M19").Value = 92.483660130719
$ws.Range("N19").Value = 47.619047619047
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 175
$ws.Range("F20").Value = 33
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 17.857142857142
$ws.Range("I20").Value = 309
$ws.Range("J20").Value = 240
$ws.Range("K20").Value = 28.75
$ws.Range("L20").Value = 164.102564102564
$ws.Range("M20").Value = 186.111111111111
$ws.Range("N20").Value = -80.480101073910
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 30.769230769230
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = 4.132231404958
$ws.Range("I21").Value = 1350
$ws.Range("J21").Value = 1338
$ws.Range("K21").Value = 0.896860986547
$ws.Range("L21").Value = 58.823529411764
$ws.Range("M21").Value = 63.834951456310
$ws.Range("N21").Value = -64.295160010579
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = 175
$ws.Range("L22").Value = 266.666666666667
$ws.Range("M22").Value = -31.25
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 16.666666666666
$ws.Range("I23").Value = 73
$ws.Range("J23").Value = 66
$ws.Range("K23").Value = 10.606060606060
$ws.Range("L23").Value = 114.705882352941
$ws.Range("M23").Value = 87.179487179487
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 28.571428571428
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = 40.697674418604
$ws.Range("I24").Value = 1005
$ws.Range("J24").Value = 1036
$ws.Range("K24").Value = -2.992277992277
$ws.Range("L24").Value = 29.677419354838
$ws.Range("M24").Value = 43.162393162393
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = -27.586206896551
$ws.Range("I25").Value = 273
$ws.Range("J25").Value = 294
$ws.Range("K25").Value = -7.142857142857
$ws.Range("L25").Value = 7.905138339920
$ws.Range("M25").Value = -21.551724137931
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -80
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = -12
$ws.Range("L26").Value = 4.761904761904
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = 0
